$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 160657
$ws.Range("C4").Value = 151680
$ws.Range("C5").Value = 8977
$ws.Range("C8").Value = 64.33
